$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.722.37'
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").Value = '2.173.94'
$ws.Range("E3").Value = '  -2.77%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.36'
$ws.Range("E5").Value = '  -1.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.53'
$ws.Range("E7").Value = '  -3.19%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  -3.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.46'
$ws.Range("E10").Value = '  -4.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").Value = '  -5.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.67'
$ws.Range("E12").Value = '  -3.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.74'
$ws.Range("E13").Value = '  -3.39%  '

$ws.Range("E14").Value = '  -3.12%  '

$ws.Range("D15").Value = '2.501.87'
$ws.Range("E15").Value = '  -2.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.41'
$ws.Range("E16").Value = '  +0.21%  '

$ws.Range("D17").Value = '2.182.16'
$ws.Range("E17").Value = '  -3.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.782'
$ws.Range("E18").Value = '  -6.95%  '

$ws.Range("D19").Value = '41.640.53'
$ws.Range("E19").Value = '  -1.14%  '

$ws.Range("E20").Value = '  -2.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.12'
$ws.Range("E21").Value = '  -3.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  -6.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.96'
$ws.Range("E23").Value = '  -13.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.39'
$ws.Range("E24").Value = '  -1.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.03'
$ws.Range("E25").Value = '  -1.00%  '

$ws.Range("E26").Value = '  +0.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.77'
$ws.Range("E27").Value = '  -5.37%  '

$ws.Range("E28").Value = '  -10.03%  '

$ws.Range("E29").Value = '  -3.57%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.89'
$ws.Range("E30").Value = '  +2.01%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("E31").Value = '  -5.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.87'
$ws.Range("E32").Value = '  -3.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.75'
$ws.Range("E33").Value = '  +9.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0777'
$ws.Range("E34").Value = '  -3.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").Value = '  -6.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.120'
$ws.Range("E36").Value = '  -3.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.34'
$ws.Range("E37").Value = '  -0.40%  '

$ws.Range("E38").Value = '  -5.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0312'
$ws.Range("E39").Value = '  +2.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.11'
$ws.Range("E40").Value = '  -9.10%  '

$ws.Range("E41").Value = '  -1.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.38'
$ws.Range("E42").Value = '  -6.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.10'
$ws.Range("E43").Value = '  -9.26%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.190'
$ws.Range("E44").Value = '  -4.93%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.47'
$ws.Range("E45").Value = '  -3.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0966'
$ws.Range("E46").Value = '  -3.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.63'
$ws.Range("E47").Value = '  -6.84%  '

$ws.Range("E48").Value = '  -4.22%  '

$ws.Range("E49").Value = '  -4.71%  '

$ws.Range("E50").Value = '  -6.78%  '

$ws.Range("E51").Value = '  -2.23%  '
